{"js": "// Fix the phase name \"Arquitectura del proyecto\" -> \"Arquitectura del programa\"\n// in the FASES table of the project-planning document.\n\n// \"Arquitectura del proyecto\" is unique in the document, so anchor on the\n// full phrase first to avoid touching the many other occurrences of the\n// word \"proyecto\" elsewhere in the text.\nconst phrase = context.document.body.search(\"Arquitectura del proyecto\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nphrase.load(\"items\");\nawait context.sync();\n\nif (phrase.items.length > 0) {\n  const target = phrase.items[0];\n\n  // Within that matched phrase, narrow down to just the word \"proyecto\"\n  // and replace it with \"programa\", leaving the rest of the text/run\n  // formatting (bold) untouched.\n  const word = target.search(\"proyecto\", { matchCase: true, matchWholeWord: false });\n  word.load(\"items\");\n  await context.sync();\n\n  if (word.items.length > 0) {\n    word.items[0].insertText(\"programa\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Fix the phase name \"Arquitectura del proyecto\" -> \"Arquitectura del programa\"\n# in the FASES table of the project-planning document.\n\n$d = $word.ActiveDocument\n\n# \"Arquitectura del proyecto\" is unique in the document, so anchor on that\n# full phrase first to avoid touching the many other occurrences of the\n# word \"proyecto\" elsewhere in the text.\n$anchor = $d.Content\n$anchor.Find.ClearFormatting()\n$anchor.Find.Text = \"Arquitectura del proyecto\"\n$found = $anchor.Find.Execute()\n\nif ($found) {\n    # Narrow down to just the word \"proyecto\" within the matched phrase and\n    # replace it with \"programa\", keeping the rest of the run/paragraph\n    # (including the bold run formatting) untouched.\n    $target = $anchor.Duplicate\n    $target.Find.ClearFormatting()\n    $target.Find.Text = \"proyecto\"\n    $target.Find.Execute() | Out-Null\n    $target.Text = \"programa\"\n}\n"}
